$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 46070

$ws.Range("C3").Value = 46070

$ws.Range("A4").Value = "A 39924-2025"
$ws.Range("B4").Value = 45891
$ws.Range("G4").Value = 1.1
$ws.Range("C4").Value = 46070

$ws.Range("A5").Value = "A 39928-2025"
$ws.Range("G5").Value = 0.9
$ws.Range("C5").Value = 46070

$ws.Range("A6").Value = "A 40001-2025"
$ws.Range("G6").Value = 0.6
$ws.Range("C6").Value = 46070

$ws.Range("A7").Value = "A 39958-2024"
$ws.Range("B7").Value = 45553
$ws.Range("G7").Value = 3.4
$ws.Range("C7").Value = 46070

$ws.Range("C8").Value = 46070

$ws.Range("A9").Value = "A 6004-2026"
$ws.Range("B9").Value = 46050
$ws.Range("G9").Value = 2.7
$ws.Range("C9").Value = 46070

$ws.Range("C10").Value = 46070

$ws.Range("A11").Value = "A 35838-2023"
$ws.Range("B11").Value = 45147
$ws.Range("G11").Value = 1.1
$ws.Range("C11").Value = 46070

$ws.Range("A12").Value = "A 2229-2023"
$ws.Range("B12").Value = 44939
$ws.Range("G12").Value = 4.3
$ws.Range("C12").Value = 46070

$ws.Range("A13").Value = "A 5528-2023"
$ws.Range("B13").Value = 44957
$ws.Range("G13").Value = 1.2
$ws.Range("C13").Value = 46070

$ws.Range("A14").Value = "A 7694-2023"
$ws.Range("B14").Value = 44967
$ws.Range("G14").Value = 2.2
$ws.Range("C14").Value = 46070

$ws.Range("A15").Value = "A 2727-2024"
$ws.Range("B15").Value = 45314
$ws.Range("G15").Value = 3.8
$ws.Range("C15").Value = 46070

$ws.Range("C16").Value = 46070

$ws.Range("C17").Value = 46070

$ws.Range("C18").Value = 46070

